$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 72 (shifts old rows 72-74 down to 73-75)
$ws.Rows(72).Insert()

# Populate the new row 72 with the fresh weekly record
$ws.Range("A72").Value = 4
$ws.Range("B72").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C72").Value = "Los Lagos"
$ws.Range("D72").Value = 44946
$ws.Range("D72").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E72").Value = 10
$ws.Range("F72").Value = "Fruta"
$ws.Range("G72").Value = 100103
$ws.Range("H72").Value = "Frutos de hueso (carozo)"
$ws.Range("I72").Value = 100103003
$ws.Range("J72").Value = "Damasco"
$ws.Range("K72").Value = "Modesto"
$ws.Range("L72").Value = "Primera"
$ws.Range("M72").Value = 600
$ws.Range("N72").Value = 20000
$ws.Range("O72").Value = 21000
$ws.Range("P72").Value = 20500
$ws.Range("Q72").Value = "$/caja 16 kilos"
$ws.Range("R72").Value = "Región de O'Higgins"
$ws.Range("S72").Value = 1281
$ws.Range("T72").Value = 16
